$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.117.99"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "1.566.79"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  +0.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.490"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0597"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "1.789.34"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").Value = "1.565.34"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "27.130.80"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "0.0₃0703"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.105"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("E30").Value = "  +5.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0472"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.86%  "
$ws.Range("D34").Value = "1.428.29"
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.61%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.807"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("E43").Value = "  +2.46%  "
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").Value = "1.708.25"
$ws.Range("E47").Value = "  +1.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("E51").Value = "  +0.27%  "
